$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the padding bug: cell B3 held a raw numeric placeholder (222)
# instead of the actual "2@2" text that its hyperlink points to. Put the
# real text in the cell so it no longer needs a separate "display" override
# on the hyperlink.
$ws.Range("B3").Value = "2@2"

# Recreate the two mailto hyperlinks on this sheet so the stale display
# text on B3's hyperlink is dropped (the cell text now matches the link).
$b2Target = "mailto:1@1"
$b3Target = "mailto:2@2"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), $b2Target)
$ws.Hyperlinks.Add($ws.Range("B3"), $b3Target)

# Recreating the hyperlinks resets the cell style, so reapply the
# "Hyperlink" cell style that both cells had before.
$ws.Range("B2").Style = "Hyperlink"
$ws.Range("B3").Style = "Hyperlink"

# --- View state: zoom in to 220% and move the selection back to B1
# (it had drifted down to E17).
$win = $wb.Windows.Item(1)
$win.Zoom = 220
$ws.Range("B1").Select() | Out-Null
